# Auto-generated Excel COM-interop edit script
# Applies the cell-value changes described by the diff to each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 196.5
$ws.Range("I33").Value = 201.66667
$ws.Range("K33").Value = 201.66667
$ws.Range("M33").Value = 27.33332999999999
$ws.Range("H51").Value = 5755.4614
$ws.Range("J51").Value = 5837.696
$ws.Range("L51").Value = 5837.696
$ws.Range("N51").Value = -6805.696
$ws.Range("H62").Value = 7112.778
$ws.Range("I62").Value = 6250.75
$ws.Range("K62").Value = 6250.75
$ws.Range("M62").Value = -5626.75
$ws.Range("H65").Value = 7112.778
$ws.Range("I65").Value = 6250.75
$ws.Range("K65").Value = 31253.75
$ws.Range("M65").Value = -28133.75
$ws.Range("H70").Value = 108015096
$ws.Range("I70").Value = 253049170
$ws.Range("J70").Value = 50001468
$ws.Range("K70").Value = 759147510
$ws.Range("L70").Value = 150004404
$ws.Range("M70").Value = -759147240
$ws.Range("N70").Value = -150004944
$ws.Range("H73").Value = 108015096
$ws.Range("I73").Value = 253049170
$ws.Range("J73").Value = 50001468
$ws.Range("K73").Value = 759147510
$ws.Range("L73").Value = 150004404
$ws.Range("M73").Value = -759146574
$ws.Range("N73").Value = -150006276
$ws.Range("H82").Value = 7622.8
$ws.Range("I82").Value = 2592.2
$ws.Range("K82").Value = 7776.599999999999
$ws.Range("M82").Value = -7370.599999999999
$ws.Range("H85").Value = 7622.8
$ws.Range("I85").Value = 2592.2
$ws.Range("K85").Value = 7776.599999999999
$ws.Range("M85").Value = -6372.599999999999
$ws.Range("H98").Value = 1595.9565
$ws.Range("I98").Value = 1271.5238
$ws.Range("J98").Value = 5002.5
$ws.Range("K98").Value = 1271.5238
$ws.Range("L98").Value = 5002.5
$ws.Range("M98").Value = 226.4762000000001
$ws.Range("N98").Value = -7998.5
$ws.Range("H122").Value = 1595.9565
$ws.Range("I122").Value = 1271.5238
$ws.Range("J122").Value = 5002.5
$ws.Range("K122").Value = 3814.5714
$ws.Range("L122").Value = 15007.5
$ws.Range("M122").Value = -1364.5714
$ws.Range("N122").Value = -19907.5
$ws.Range("H131").Value = 2071245.6
$ws.Range("I131").Value = 2533.2856
$ws.Range("K131").Value = 7599.8568
$ws.Range("M131").Value = -2559.8568
$ws.Range("H137").Value = 11366352
$ws.Range("I137").Value = 15627201
$ws.Range("J137").Value = 4090.0833
$ws.Range("K137").Value = 46881603
$ws.Range("L137").Value = 12270.2499
$ws.Range("M137").Value = -46879053
$ws.Range("N137").Value = -17370.2499
$ws.Range("H138").Value = 6822.4443
$ws.Range("J138").Value = 9610.105
$ws.Range("L138").Value = 28830.315
$ws.Range("N138").Value = -39110.315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5012.467
$ws.Range("I2").Value = 4698.8184
$ws.Range("J2").Value = 5875
$ws.Range("K2").Value = 4698.8184
$ws.Range("L2").Value = 5875
$ws.Range("M2").Value = -4585.8184
$ws.Range("N2").Value = -6101
$ws.Range("H30").Value = 39999
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = 9
$ws.Range("M30").Value = 141
$ws.Range("H32").Value = 4173.1113
$ws.Range("I32").Value = 3950.7932
$ws.Range("K32").Value = 3950.7932
$ws.Range("M32").Value = -3663.7932
$ws.Range("H61").Value = 78891110
$ws.Range("I61").Value = 87502504
$ws.Range("K61").Value = 87502504
$ws.Range("M61").Value = -87502292
$ws.Range("H97").Value = 1201.4
$ws.Range("I97").Value = 1040.4762
$ws.Range("J97").Value = 1576.8889
$ws.Range("K97").Value = 1040.4762
$ws.Range("L97").Value = 1576.8889
$ws.Range("M97").Value = -544.4762000000001
$ws.Range("N97").Value = -2568.8889
$ws.Range("H116").Value = 5012.467
$ws.Range("I116").Value = 4698.8184
$ws.Range("J116").Value = 5875
$ws.Range("K116").Value = 4698.8184
$ws.Range("L116").Value = 5875
$ws.Range("M116").Value = -2404.8184
$ws.Range("N116").Value = -10463
$ws.Range("H136").Value = 78891110
$ws.Range("I136").Value = 87502504
$ws.Range("K136").Value = 262507512
$ws.Range("M136").Value = -262504962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5012.467
$ws.Range("I3").Value = 4698.8184
$ws.Range("J3").Value = 5875
$ws.Range("K3").Value = 4698.8184
$ws.Range("L3").Value = 5875
$ws.Range("M3").Value = -4584.8184
$ws.Range("N3").Value = -6103
$ws.Range("H114").Value = 74997
$ws.Range("J114").Value = 74997
$ws.Range("L114").Value = 74997
$ws.Range("N114").Value = -83675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4339.357
$ws.Range("I62").Value = 4128.5454
$ws.Range("J62").Value = 5112.3335
$ws.Range("K62").Value = 4128.5454
$ws.Range("L62").Value = 5112.3335
$ws.Range("M62").Value = -3504.5454
$ws.Range("N62").Value = -6360.3335
$ws.Range("H65").Value = 4339.357
$ws.Range("I65").Value = 4128.5454
$ws.Range("J65").Value = 5112.3335
$ws.Range("K65").Value = 20642.727
$ws.Range("L65").Value = 25561.6675
$ws.Range("M65").Value = -17522.727
$ws.Range("N65").Value = -31801.6675
$ws.Range("H132").Value = 2100
$ws.Range("I132").Value = 2100
$ws.Range("K132").Value = 6300
$ws.Range("M132").Value = -3770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2419.3
$ws.Range("J68").Value = 2554.889
$ws.Range("L68").Value = 7664.667
$ws.Range("N68").Value = -9286.667000000001
$ws.Range("H71").Value = 2419.3
$ws.Range("J71").Value = 2554.889
$ws.Range("L71").Value = 22994.001
$ws.Range("N71").Value = -31106.001
$ws.Range("H81").Value = 11261
$ws.Range("I81").Value = 225
$ws.Range("J81").Value = 33333
$ws.Range("K81").Value = 675
$ws.Range("L81").Value = 99999
$ws.Range("M81").Value = 448
$ws.Range("N81").Value = -102245
$ws.Range("H84").Value = 11261
$ws.Range("I84").Value = 225
$ws.Range("J84").Value = 33333
$ws.Range("K84").Value = 2025
$ws.Range("L84").Value = 299997
$ws.Range("M84").Value = 3591
$ws.Range("N84").Value = -311229
$ws.Range("H97").Value = 1740.6666
$ws.Range("J97").Value = 1559.7142
$ws.Range("L97").Value = 4679.142599999999
$ws.Range("N97").Value = -5671.142599999999
$ws.Range("H107").Value = 3795598.2
$ws.Range("I107").Value = 2023
$ws.Range("K107").Value = 6069
$ws.Range("M107").Value = -4149
$ws.Range("H109").Value = 5557.1113
$ws.Range("I109").Value = 954.4286
$ws.Range("K109").Value = 2863.2858
$ws.Range("M109").Value = -1823.2858
$ws.Range("H117").Value = 5264.1665
$ws.Range("I117").Value = 2000
$ws.Range("J117").Value = 5560.909
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 16682.727
$ws.Range("M117").Value = -2558
$ws.Range("N117").Value = -23566.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2925.7856
$ws.Range("I122").Value = 2774.7778
$ws.Range("K122").Value = 8324.3334
$ws.Range("M122").Value = -5874.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2780216.5
$ws.Range("I68").Value = 4631762
$ws.Range("J68").Value = 2898.1667
$ws.Range("K68").Value = 4631762
$ws.Range("L68").Value = 2898.1667
$ws.Range("M68").Value = -4631013
$ws.Range("N68").Value = -4396.1667
$ws.Range("H71").Value = 2780216.5
$ws.Range("I71").Value = 4631762
$ws.Range("J71").Value = 2898.1667
$ws.Range("K71").Value = 23158810
$ws.Range("L71").Value = 14490.8335
$ws.Range("M71").Value = -23155066
$ws.Range("N71").Value = -21978.8335
$ws.Range("H93").Value = 1427206.9
$ws.Range("I93").Value = 1003.11536
$ws.Range("K93").Value = 1003.11536
$ws.Range("M93").Value = 244.88464
$ws.Range("H132").Value = 5593.3
$ws.Range("J132").Value = 5603.778
$ws.Range("L132").Value = 16811.334
$ws.Range("N132").Value = -21871.334
$ws.Range("H136").Value = 3355.1292
$ws.Range("I136").Value = 3216.6667
$ws.Range("J136").Value = 3645.9
$ws.Range("K136").Value = 9650.000100000001
$ws.Range("L136").Value = 10937.7
$ws.Range("M136").Value = -7100.000100000001
$ws.Range("N136").Value = -16037.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3026.32
$ws.Range("I132").Value = 2483.3125
$ws.Range("J132").Value = 3991.6667
$ws.Range("K132").Value = 7449.9375
$ws.Range("L132").Value = 11975.0001
$ws.Range("M132").Value = -4919.9375
$ws.Range("N132").Value = -17035.0001
